$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 482, shifting existing rows 482:560 down to 483:561.
$ws.Range("A482:T482").Insert()

# Populate the newly inserted row 482 with the new weekly record.
$ws.Range("A482").Value = 9
$ws.Range("B482").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C482").Value = "Metropolitana"
$ws.Range("D482").Value = 44491
$ws.Range("E482").Value = 13
$ws.Range("F482").Value = "Fruta"
$ws.Range("G482").Value = 100104
$ws.Range("H482").Value = "Frutos de pepita"
$ws.Range("I482").Value = 100104005
$ws.Range("J482").Value = "Pera"
$ws.Range("K482").Value = "Packham's Triumph"
$ws.Range("L482").Value = "Calibre 80"
$ws.Range("M482").Value = 300
$ws.Range("N482").Value = 19000
$ws.Range("O482").Value = 19000
$ws.Range("P482").Value = 19000
$ws.Range("Q482").Value = "`$/caja 18 kilos embalada"
$ws.Range("R482").Value = "Provincia de Curicó"
$ws.Range("S482").Value = 1056
$ws.Range("T482").Value = 18

Write-Output "Row inserted and populated"
